# Update DLM status to Co-PI
#
# The "grants" sheet gets a new current entry reflecting a role change to
# Co-Principal Investigator on the Dynamic Learning Maps Alternate
# Assessment System contract. The new entry is inserted as the first data
# row (row 2), pushing the existing three grant rows down one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grants")

# Insert a new blank row above the current row 2 (shifts rows 2-4 -> 3-5).
$ws.Rows.Item(2).Insert()

# Populate the new row with the Co-PI / DLM entry.
$ws.Cells.Item(2, 1).Value = "Currently Funded"
$ws.Cells.Item(2, 2).Value = "Co-Principal Investigator"
$ws.Cells.Item(2, 3).Value = "Dynamic Learning Maps Alternate Assessment System"
$ws.Cells.Item(2, 4).Value = "Ongoing state contracts"
$ws.Cells.Item(2, 8).Value = "Meagan Karvonen"

# New row is a single-line row (unlike the multi-line ht=34 rows below it).
$ws.Rows.Item(2).RowHeight = 17

# Make the grants sheet the active tab/sheet, with the previously-added
# row selected, matching the saved workbook view state after the edit.
$ws.Activate()
$ws.Range("D3").Select()
